# Updated symbol list on Sat Dec 24 14:26:54 UTC 2022 with GitHub Actions
#
# This script reproduces the price/ranking refresh that the GitHub Actions
# job performs on the cryptos worksheet: the "Price" column (D) is
# refreshed with newly scraped numeric-looking values (kept as TEXT, same
# as the source data), and several rows shift which coin/link/volume they
# describe because the underlying ranking order changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking price into column D while preserving it
# as TEXT (matching the workbook's existing inline-string / text storage)
# instead of letting Excel auto-convert it to a Number.
function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Simple price (column D) refreshes ---
Set-TextValue "D2" "244.55"
Set-TextValue "D4" "5.383"
Set-TextValue "D6" "3.392"
Set-TextValue "D7" "0.8142"
Set-TextValue "D8" "0.9606"

# --- Rows 9-17: coin ranking reshuffled, each row now describes the coin
#     that used to be one row below it (plus fresh price data) ---
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1424"
$ws.Range("E9").Value = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07417"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.03319"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03044"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09409"
$ws.Range("E13").Value = "12BitMartTokenBMX"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "4.000"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001588"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04816"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D17" "0.0005902"
$ws.Range("E17").Value = "16OneONE"

# --- More simple price (column D) refreshes ---
Set-TextValue "D18" "0.005462"
Set-TextValue "D19" "0.004147"
Set-TextValue "D20" "0.0009889"

# --- Row 27: "Best in 24h" badge moved onto UpBots ---
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

Set-TextValue "D40" "0.04034"

# --- Rows 41-43: another coin ranking reshuffle ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1074"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002720"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003056"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Final simple price (column D) refreshes ---
Set-TextValue "D44" "0.005684"
Set-TextValue "D45" "0.00005217"
Set-TextValue "D47" "0.8603"
Set-TextValue "D48" "0.007198"
